# 11-01-2018 work done mostly on python
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "python basics" entry (row 35) with the expanded description
$ws.Range("B35").Value = "python basics(overview,basic syntax,variables,operators,decision making,loops"

# Row 36 gets a new, different python-basics follow-up entry
$ws.Range("B36").Value = "python basics( number string list tuple dictionary d&t"

# New row 37 for the next day of work
$ws.Range("A37").Value = 43111
$ws.Range("B37").Value = "python: functions modules file I/O exceptions,OOP reg"

# Give the new B37 cell the same centered/wrapped formatting as the rows above it
$ws.Range("B36").Copy() | Out-Null
$ws.Range("B37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row heights for the wrapped-text entries
$ws.Rows.Item(35).RowHeight = 45
$ws.Rows.Item(36).RowHeight = 30
$ws.Rows.Item(37).RowHeight = 30

# Update the view to reflect scrolling down to the new rows
$excel.ActiveWindow.ScrollRow = 30
$ws.Range("D36").Select()
